$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that look like pure numbers need a leading apostrophe so Excel
# keeps them stored as text (matching the original shared-string cells);
# values that already contain non-numeric characters (e.g. "***", "*")
# are kept as text automatically.
$ws.Range("B2").Value = "'0.17"
$ws.Range("B3").Value = "'-0.01"
$ws.Range("B4").Value = "'-0.09"
$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "'0.98"
$ws.Range("D2").Value = "'-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

# Reset style on the apostrophe-entered cells so Excel doesn't leave a
# "quote prefix" cell format applied (keeps styling identical to source).
$ws.Range("B2:D4").Style = "Normal"
